$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2452590214026884
$ws.Range("C2").Value = 0.03524592337566901
$ws.Range("E2").Value = 0.1649013415002969
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002425384863620384
$ws.Range("I2").Value = 0.4984365121658598
$ws.Range("K2").Value = 0.246392164129702
$ws.Range("M2").Value = 0.2174129915600531
$ws.Range("N2").Value = 1.38848349076085
$ws.Range("O2").Value = 2.209038766487879
$ws.Range("B3").Value = 0.2161548119744907
$ws.Range("C3").Value = 0.03128430103414814
$ws.Range("E3").Value = 0.1534863313712123
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002427427657781619
$ws.Range("I3").Value = 0.5036256256496543
$ws.Range("K3").Value = 0.2155334280826935
$ws.Range("M3").Value = 0.195872510153059
$ws.Range("N3").Value = 1.403396919388125
$ws.Range("O3").Value = 2.225653602028444
$ws.Range("B4").Value = 0.1982719915989719
$ws.Range("C4").Value = 0.02883719686003872
$ws.Range("E4").Value = 0.1465879726479926
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002428748654920306
$ws.Range("I4").Value = 0.5070795527495093
$ws.Range("K4").Value = 0.1965487396685575
$ws.Range("M4").Value = 0.1827167963917304
$ws.Range("N4").Value = 1.413043558532442
$ws.Range("O4").Value = 2.237138901456646
$ws.Range("B5").Value = 0.190981833954595
$ws.Range("C5").Value = 0.02783632598183772
$ws.Range("E5").Value = 0.1438044538709846
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002429303796654356
$ws.Range("I5").Value = 0.5085543681700706
$ws.Range("K5").Value = 0.1888033548759722
$ws.Range("M5").Value = 0.1773734123769017
$ws.Range("N5").Value = 1.417097800820738
$ws.Range("O5").Value = 2.242141902143331
$ws.Range("B6").Value = 0.189771154754709
$ws.Range("C6").Value = 0.02766991215584369
$ws.Range("E6").Value = 0.1433439158568603
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002429396994860515
$ws.Range("I6").Value = 0.5088033248746555
$ws.Range("K6").Value = 0.1875167099149451
$ws.Range("M6").Value = 0.1764872156420694
$ws.Range("N6").Value = 1.417778443500927
$ws.Range("O6").Value = 2.242992129234381
$ws.Range("B7").Value = 0.1981736845949627
$ws.Range("C7").Value = 0.02882371353041435
$ws.Range("E7").Value = 0.1465503215371342
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002428756073390672
$ws.Range("I7").Value = 0.5070991701166747
$ws.Range("K7").Value = 0.1964443184781004
$ws.Range("M7").Value = 0.1826446620536828
$ws.Range("N7").Value = 1.413097736782763
$ws.Range("O7").Value = 2.237205067491146
$ws.Range("B8").Value = 0.2352267816405345
$ws.Range("C8").Value = 0.03388301286470607
$ws.Range("E8").Value = 0.1609424059019204
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002426075402202697
$ws.Range("I8").Value = 0.5001701230989859
$ws.Range("K8").Value = 0.2357600496829662
$ws.Range("M8").Value = 0.2099712607768112
$ws.Range("N8").Value = 1.393524004458435
$ws.Range("O8").Value = 2.214501030320775
$ws.Range("B9").Value = 0.3077712575458804
$ws.Range("C9").Value = 0.04368723574260969
$ws.Range("E9").Value = 0.1900511884586251
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002421345706616274
$ws.Range("I9").Value = 0.4887082553233668
$ws.Range("K9").Value = 0.312547795218336
$ws.Range("M9").Value = 0.2641179795588755
$ws.Range("N9").Value = 1.359024791015553
$ws.Range("O9").Value = 2.180172234568047
$ws.Range("B10").Value = 0.3609827820231999
$ws.Range("C10").Value = 0.05081865061498547
$ws.Range("E10").Value = 0.2119924022077484
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002418188971907411
$ws.Range("I10").Value = 0.4815849029922674
$ws.Range("K10").Value = 0.3687607201642322
$ws.Range("M10").Value = 0.3042475899942971
$ws.Range("N10").Value = 1.336044397771296
$ws.Range("O10").Value = 2.161177758619303
$ws.Range("B11").Value = 0.3851681545645818
$ws.Range("C11").Value = 0.0540472942710295
$ws.Range("E11").Value = 0.222097864257961
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002416821308223405
$ws.Range("I11").Value = 0.4786264481774012
$ws.Range("K11").Value = 0.3942867262827292
$ws.Range("M11").Value = 0.322580898090699
$ws.Range("N11").Value = 1.32610338280762
$ws.Range("O11").Value = 2.153891614035018
$ws.Range("B12").Value = 0.3943231567109251
$ws.Range("C12").Value = 0.05526765322539973
$ws.Range("E12").Value = 0.2259426301113052
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002416313187645935
$ws.Range("I12").Value = 0.4775467431072471
$ws.Range("K12").Value = 0.4039458668828217
$ws.Range("M12").Value = 0.3295345275117469
$ws.Range("N12").Value = 1.322412729903359
$ws.Range("O12").Value = 2.151327539745836
$ws.Range("B13").Value = 0.3923516253242667
$ws.Range("C13").Value = 0.05500492806031332
$ws.Range("E13").Value = 0.2251137856067658
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002416422185974154
$ws.Range("I13").Value = 0.4777774705729421
$ws.Range("K13").Value = 0.401865916891694
$ws.Range("M13").Value = 0.3280364412696173
$ws.Range("N13").Value = 1.323204295348486
$ws.Range("O13").Value = 2.151871081010938
$ws.Range("B14").Value = 0.3859214144665088
$ws.Range("C14").Value = 0.05414773939031647
$ws.Range("E14").Value = 0.2224138127847723
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002416779309048385
$ws.Range("I14").Value = 0.4785368062853301
$ws.Range("K14").Value = 0.3950815325742099
$ws.Range("M14").Value = 0.3231527530303353
$ws.Range("N14").Value = 1.32579827138612
$ws.Range("O14").Value = 2.153676755967865
$ws.Range("B15").Value = 0.381982258042143
$ws.Range("C15").Value = 0.0536223914865559
$ws.Range("E15").Value = 0.2207623568199253
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002416999329844389
$ws.Range("I15").Value = 0.4790072097371798
$ws.Range("K15").Value = 0.3909249756337942
$ws.Range("M15").Value = 0.3201628109778341
$ws.Range("N15").Value = 1.327396768019756
$ws.Range("O15").Value = 2.154808190157695
$ws.Range("B16").Value = 0.3594017511158825
$ws.Range("C16").Value = 0.05060733695405872
$ws.Range("E16").Value = 0.2113345024415665
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002418279723463415
$ws.Range("I16").Value = 0.4817839234414834
$ws.Range("K16").Value = 0.3670915803625689
$ws.Range("M16").Value = 0.3030510348344961
$ws.Range("N16").Value = 1.336704390492113
$ws.Range("O16").Value = 2.161681203408492
$ws.Range("B17").Value = 0.3455436681520609
$ws.Range("C17").Value = 0.04875371196294509
$ws.Range("E17").Value = 0.2055827773338663
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002419082677887509
$ws.Range("I17").Value = 0.4835596063058496
$ws.Range("K17").Value = 0.3524585771471607
$ws.Range("M17").Value = 0.2925735059164936
$ws.Range("N17").Value = 1.342545705659372
$ws.Range("O17").Value = 2.166244681692973
$ws.Range("B18").Value = 0.3375709346721578
$ws.Range("C18").Value = 0.04768610152089536
$ws.Range("E18").Value = 0.2022862259807354
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002419550952074267
$ws.Range("I18").Value = 0.4846074718956537
$ws.Range("K18").Value = 0.3440378019716661
$ws.Range("M18").Value = 0.2865544660622987
$ws.Range("N18").Value = 1.345953747767194
$ws.Range("O18").Value = 2.168996954278029
$ws.Range("B19").Value = 0.3348711871402372
$ws.Range("C19").Value = 0.04732437819534141
$ws.Range("E19").Value = 0.2011720726049262
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002419710608629222
$ws.Range("I19").Value = 0.4849668179756002
$ws.Range("K19").Value = 0.3411859546368987
$ws.Range("M19").Value = 0.284517787291243
$ws.Range("N19").Value = 1.347115941814575
$ws.Range("O19").Value = 2.169950713793853
$ws.Range("B20").Value = 0.347019088037996
$ws.Range("C20").Value = 0.04895118427185707
$ws.Range("E20").Value = 0.2061938472868832
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002418996536315888
$ws.Range("I20").Value = 0.4833678347248309
$ws.Range("K20").Value = 0.3540167293880927
$ws.Range("M20").Value = 0.293688095510177
$ws.Range("N20").Value = 1.341918891406564
$ws.Range("O20").Value = 2.165745696066224
$ws.Range("B21").Value = 0.3878102213622299
$ws.Range("C21").Value = 0.05439957803729101
$ws.Range("E21").Value = 0.223206369089894
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002416674148194182
$ws.Range("I21").Value = 0.4783126687974111
$ws.Range("K21").Value = 0.3970744644224453
$ws.Range("M21").Value = 0.3245869068630114
$ws.Range("N21").Value = 1.325034354656875
$ws.Range("O21").Value = 2.153141090400908
$ws.Range("B22").Value = 0.4144491158656933
$ws.Range("C22").Value = 0.05794723959131431
$ws.Range("E22").Value = 0.2344303055570123
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002415213341679049
$ws.Range("I22").Value = 0.4752454812466276
$ws.Range("K22").Value = 0.4251741475201243
$ws.Range("M22").Value = 0.3448463526211256
$ws.Range("N22").Value = 1.314429481399657
$ws.Range("O22").Value = 2.146040154974031
$ws.Range("B23").Value = 0.4002334680519652
$ws.Range("C23").Value = 0.05605500272224617
$ws.Range("E23").Value = 0.2284301900691048
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002415987800053049
$ws.Range("I23").Value = 0.4768608273516257
$ws.Range("K23").Value = 0.4101807232441672
$ws.Range("M23").Value = 0.3340275428816639
$ws.Range("N23").Value = 1.320050128757551
$ws.Range("O23").Value = 2.149725948680612
$ws.Range("B24").Value = 0.346352067918815
$ws.Range("C24").Value = 0.04886191307105037
$ws.Range("E24").Value = 0.2059175508294118
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.00241903546020985
$ws.Range("I24").Value = 0.4834544505506244
$ws.Range("K24").Value = 0.3533123138526832
$ws.Range("M24").Value = 0.2931841751294186
$ws.Range("N24").Value = 1.342202119031207
$ws.Range("O24").Value = 2.16597088669117
$ws.Range("B25").Value = 0.2881601709947859
$ws.Range("C25").Value = 0.04104750517937816
$ws.Range("E25").Value = 0.1820801033752772
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002422569113581719
$ws.Range("I25").Value = 0.4915812262903465
$ws.Range("K25").Value = 0.2918093285035468
$ws.Range("M25").Value = 0.2494092467197291
$ws.Range("N25").Value = 1.36794207457628
$ws.Range("O25").Value = 2.188366481341049
